# Add a new "company" worksheet at the end of the workbook, matching the
# "adding test case for add new company" commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "company"

$headerRow = New-Object 'object[,]' 1,24
$row2 = New-Object 'object[,]' 1,24
$row3 = New-Object 'object[,]' 1,24
$headerRow[0,0] = 'Company'
$headerRow[0,1] = 'Industry'
$headerRow[0,2] = 'AnnualRevenue'
$headerRow[0,3] = 'Employees'
$headerRow[0,4] = 'Status'
$headerRow[0,5] = 'Category'
$headerRow[0,6] = 'Priority'
$headerRow[0,7] = 'Source'
$headerRow[0,8] = 'Identifier'
$headerRow[0,9] = 'VAT/TaxNumber'
$headerRow[0,10] = 'AddressTitle'
$headerRow[0,11] = 'DefaultAddress'
$headerRow[0,12] = 'City'
$headerRow[0,13] = 'State'
$headerRow[0,14] = 'Zip'
$headerRow[0,15] = 'Country'
$headerRow[0,16] = 'Tags '
$headerRow[0,17] = 'Description'
$headerRow[0,18] = 'Phone'
$headerRow[0,19] = 'Fax'
$headerRow[0,20] = 'Website'
$headerRow[0,21] = 'Email'
$headerRow[0,22] = 'Symbol'
$headerRow[0,23] = 'ParentCompany'
$row2[0,0] = 'Amazon'
$row2[0,1] = 'E-Commerce'
$row2[0,2] = 100000
$row2[0,3] = 5000
$row2[0,4] = 'Active'
$row2[0,5] = 'Client'
$row2[0,6] = 'High'
$row2[0,7] = 'Email'
$row2[0,8] = 'test identifier -1'
$row2[0,9] = 'Vat-1'
$row2[0,10] = 'Main Office'
$row2[0,11] = 'Add -1, Add-2, Add-3, Add-4'
$row2[0,12] = 'san diego'
$row2[0,13] = 'Texas'
$row2[0,14] = 121212
$row2[0,15] = 'US'
$row2[0,16] = 'Tag 1'
$row2[0,17] = 'Test desc 1'
$row2[0,18] = 1125658911
$row2[0,19] = 798645132
$row2[0,20] = 'www.amazon.com1'
$row2[0,21] = 'test@test.com'
$row2[0,22] = 'symbol -1'
$row2[0,23] = 'parent one'
$row3[0,0] = 'Microsoft'
$row3[0,1] = 'IT'
$row3[0,2] = 9000000
$row3[0,3] = 8000
$row3[0,4] = 'Hot'
$row3[0,5] = 'Partner'
$row3[0,6] = 'Medium'
$row3[0,7] = 'Internet'
$row3[0,8] = 'test identifier -1'
$row3[0,9] = 'Vat-2'
$row3[0,10] = 'Branch office'
$row3[0,11] = 'Add -7, Add-8, Add-9, Add-10'
$row3[0,12] = 'los angeles'
$row3[0,13] = 'New York'
$row3[0,14] = 323232
$row3[0,15] = 'UK'
$row3[0,16] = 'Tag 2'
$row3[0,17] = 'Test desc 2'
$row3[0,18] = 9883465555
$row3[0,19] = 364665132
$row3[0,20] = 'www.microsoft.com1'
$row3[0,21] = 'test@tes1.com'
$row3[0,22] = 'symbol -2'
$row3[0,23] = 'parent two'

$ws.Range("A1:X1").Value2 = $headerRow
$ws.Range("A2:X2").Value2 = $row2
$ws.Range("A3:X3").Value2 = $row3

# Website / Email columns carry live hyperlinks (U = Website, V = Email)
$ws.Hyperlinks.Add($ws.Range("U2"), "http://www.amazon.com1")
$ws.Hyperlinks.Add($ws.Range("U3"), "http://www.microsoft.com1")
$ws.Hyperlinks.Add($ws.Range("V2"), "mailto:test@test.com")
$ws.Hyperlinks.Add($ws.Range("V3"), "mailto:test@tes1.com")

# Re-apply the built-in Hyperlink style so the cells reuse the workbook's
# existing style slot instead of leaving a fresh (duplicate) one selected.
$ws.Range("U2:V3").Style = "Hyperlink"

# Column sizing to best-fit the new data (matches the "bestFit" columns
# Excel computes automatically when a sheet like this is authored).
$ws.Columns.Item("A:X").AutoFit()

# Scroll/selection state recorded for the new tab.
$ws.Range("X8").Select()

